# Update the LR-pairs export with the new TPM-based values, including the
# addition of "Resolving-Mac" as both a sending and target cluster, which
# expands the table from 12 to 16 data rows (plus header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("ECs", "Clec11a", "Itgb1", "ECs", 2, 0.6666666666666666, 0.152918, 0.458754, 0.007158013163202275, 0.007158013163202275, 3, 1, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843, 22.28026647293733, 200.522398256436, 0.002051442014017593, 0.002051442014017593)
  ,@("ECs", "Clec11a", "Itgb1", "FAPs", 2, 0.6666666666666666, 0.152918, 0.458754, 0.007158013163202275, 0.007158013163202275, 3, 1, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841, 25.81251293238134, 232.312616391432, 0.002376671463116413, 0.002376671463116413)
  ,@("ECs", "Clec11a", "Itgb1", "MuSCs", 2, 0.6666666666666666, 0.152918, 0.458754, 0.007158013163202275, 0.007158013163202275, 3, 1, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105, 19.59279531931734, 176.335157873856, 0.001803994738524143, 0.001803994738524143)
  ,@("ECs", "Clec11a", "Itgb1", "Resolving-Mac", 2, 0.6666666666666666, 0.152918, 0.458754, 0.007158013163202275, 0.007158013163202275, 3, 1, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212, 10.056052678522, 90.504474106698, 0.0009259049475441256, 0.0009259049475441256)
  ,@("FAPs", "Clec11a", "Itgb1", "ECs", 3, 1, 20.82581466666667, 62.47744400000001, 0.9748457050079848, 0.9748457050079848, 3, 1, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843, 3034.336705223322, 27309.0303470099, 0.2793847106510928, 0.2793847106510928)
  ,@("FAPs", "Clec11a", "Itgb1", "FAPs", 3, 1, 20.82581466666667, 62.47744400000001, 0.9748457050079848, 0.9748457050079848, 3, 1, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841, 3515.39132352444, 31638.52191171996, 0.3236775226880938, 0.3236775226880938)
  ,@("FAPs", "Clec11a", "Itgb1", "MuSCs", 3, 1, 20.82581466666667, 62.47744400000001, 0.9748457050079848, 0.9748457050079848, 3, 1, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105, 2668.331551040669, 24014.98395936602, 0.2456850082014257, 0.2456850082014257)
  ,@("FAPs", "Clec11a", "Itgb1", "Resolving-Mac", 3, 1, 20.82581466666667, 62.47744400000001, 0.9748457050079848, 0.9748457050079848, 3, 1, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212, 1369.528043534026, 12325.75239180623, 0.1260984634673726, 0.1260984634673726)
  ,@("MuSCs", "Clec11a", "Itgb1", "ECs", 3, 1, 0.200477, 0.601431, 0.009384225564807953, 0.009384225564807953, 3, 1, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843, 29.20964818853934, 262.8868336968541, 0.002689460630169144, 0.002689460630169144)
  ,@("MuSCs", "Clec11a", "Itgb1", "FAPs", 3, 1, 0.200477, 0.601431, 0.009384225564807953, 0.009384225564807953, 3, 1, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841, 33.84045799150534, 304.564121923548, 0.003115839632425151, 0.003115839632425151)
  ,@("MuSCs", "Clec11a", "Itgb1", "MuSCs", 3, 1, 0.200477, 0.601431, 0.009384225564807953, 0.009384225564807953, 3, 1, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105, 25.68634710910934, 231.177123981984, 0.002365054821506329, 0.002365054821506329)
  ,@("MuSCs", "Clec11a", "Itgb1", "Resolving-Mac", 3, 1, 0.200477, 0.601431, 0.009384225564807953, 0.009384225564807953, 3, 1, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212, 13.183583834683, 118.652254512147, 0.001213870480707331, 0.001213870480707331)
  ,@("Resolving-Mac", "Clec11a", "Itgb1", "ECs", 2, 0.6666666666666666, 0.183981, 0.5519430000000001, 0.008612056264005009, 0.008612056264005009, 3, 1, 145.7007446666667, 437.1022340000001, 0.2865937750105843, 0.2865937750105843, 26.80616870451801, 241.2555183406621, 0.002468161715304744, 0.002468161715304744)
  ,@("Resolving-Mac", "Clec11a", "Itgb1", "FAPs", 2, 0.6666666666666666, 0.183981, 0.5519430000000001, 0.008612056264005009, 0.008612056264005009, 3, 1, 168.7997026666667, 506.3991080000001, 0.3320294904365841, 0.3320294904365841, 31.05593809631601, 279.5034428668441, 0.002859456652948776, 0.002859456652948776)
  ,@("Resolving-Mac", "Clec11a", "Itgb1", "MuSCs", 2, 0.6666666666666666, 0.183981, 0.5519430000000001, 0.008612056264005009, 0.008612056264005009, 3, 1, 128.1261546666667, 384.378464, 0.2520245069956105, 0.2520245069956105, 23.57277806172801, 212.155002555552, 0.002170449234154321, 0.002170449234154321)
  ,@("Resolving-Mac", "Clec11a", "Itgb1", "Resolving-Mac", 2, 0.6666666666666666, 0.183981, 0.5519430000000001, 0.008612056264005009, 0.008612056264005009, 3, 1, 65.761079, 197.283237, 0.1293522275572212, 0.1293522275572212, 12.098789075499, 108.889101679491, 0.001113988661597168, 0.001113988661597168)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $rowNum = $i + 2
  $rowVals = $data[$i]
  for ($j = 0; $j -lt $rowVals.Count; $j++) {
    $colNum = $j + 1
    $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
  }
}

